$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells that look like plain numbers to remain as text,
# matching the source data (which stores Price/Volume as strings).
$ws.Cells.Item(2, 4).Value = '60.404.98'
$ws.Cells.Item(2, 5).Value = '  +0.02%  '
$ws.Cells.Item(3, 4).Value = '2.338.91'
$ws.Cells.Item(3, 5).Value = '  -0.39%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '548.62'
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '131.43'
$ws.Cells.Item(6, 5).Value = '  -0.83%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '1.00'
$ws.Cells.Item(7, 5).Value = '  +0.01%  '
$ws.Cells.Item(8, 5).Value = '  -1.14%  '
$ws.Cells.Item(9, 4).Value = '2.337.60'
$ws.Cells.Item(9, 5).Value = '  -0.33%  '
$ws.Cells.Item(10, 5).Value = '  +0.91%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '5.62'
$ws.Cells.Item(11, 5).Value = '  +1.63%  '
$ws.Cells.Item(12, 5).Value = '  -0.70%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.338'
$ws.Cells.Item(13, 5).Value = '  +0.69%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '23.74'
$ws.Cells.Item(14, 5).Value = '  -0.88%  '
$ws.Cells.Item(15, 4).Value = '2.755.25'
$ws.Cells.Item(15, 5).Value = '  -0.31%  '
$ws.Cells.Item(16, 4).Value = '60.358.48'
$ws.Cells.Item(17, 5).Value = '  +0.85%  '
$ws.Cells.Item(18, 4).Value = '2.341.57'
$ws.Cells.Item(18, 5).Value = '  +0.05%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '10.70'
$ws.Cells.Item(19, 5).Value = '  +0.37%  '
$ws.Cells.Item(20, 5).Value = '  -1.39%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '315.29'
$ws.Cells.Item(21, 5).Value = '  +0.09%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.60'
$ws.Cells.Item(22, 5).Value = '  -3.27%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '64.26'
$ws.Cells.Item(24, 5).Value = '  +1.03%  '
$ws.Cells.Item(25, 5).Value = '  -1.54%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '0.998'
$ws.Cells.Item(26, 5).Value = '  -0.10%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '7.98'
$ws.Cells.Item(27, 5).Value = '  +1.01%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '1.39'
$ws.Cells.Item(28, 5).Value = '  +1.42%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '1.27'
$ws.Cells.Item(29, 5).Value = '  +8.76%  '
$ws.Cells.Item(30, 5).Value = '  -0.54%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '171.22'
$ws.Cells.Item(31, 5).Value = '  -0.26%  '
$ws.Cells.Item(32, 4).Value = '0.0₃0734'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '6.13'
$ws.Cells.Item(33, 5).Value = '  +3.04%  '
$ws.Cells.Item(34, 5).Value = '  -2.87%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.385'
$ws.Cells.Item(35, 5).Value = '  +0.77%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '18.12'
$ws.Cells.Item(36, 5).Value = '  +0.45%  '
$ws.Cells.Item(38, 5).Value = '  +0.01%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '4.11'
$ws.Cells.Item(39, 5).Value = '  -1.59%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '324.96'
$ws.Cells.Item(40, 5).Value = '  +0.00%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '38.18'
$ws.Cells.Item(41, 5).Value = '  +0.33%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '1.54'
$ws.Cells.Item(42, 5).Value = '  +0.44%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '137.90'
$ws.Cells.Item(43, 5).Value = '  -2.59%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '3.52'
$ws.Cells.Item(44, 5).Value = '  +1.67%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.0951'
$ws.Cells.Item(45, 5).Value = '  -0.02%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '19.40'
$ws.Cells.Item(46, 5).Value = '  -0.91%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.568'
$ws.Cells.Item(47, 5).Value = '  +1.07%  '
$ws.Cells.Item(48, 5).Value = '  +0.26%  '
$ws.Cells.Item(49, 5).Value = '  +1.25%  '
$ws.Cells.Item(50, 5).Value = '  +4.54%  '
$ws.Cells.Item(51, 5).Value = '  -0.86%  '
